$wb = $excel.ActiveWorkbook

$wsSteps = $wb.Worksheets.Item("Test Steps")
$wsCases = $wb.Worksheets.Item("Test Cases")

# Update selection on "Test Steps" sheet to E2
[void]$wsSteps.Range("E2").Select()

# Fill column G (rows 2-9) on "Test Steps" with PASS
$wsSteps.Range("G2:G9").Value = "PASS"
$wsSteps.Range("G2:G9").Style = "Normal"

# Fill D2 on "Test Cases" with PASS
$wsCases.Range("D2").Value = "PASS"
$wsCases.Range("D2").Style = "Normal"
